# Revise the quiz data file to be "R-friendly": lower-case the
# Female/Male category labels (stored as shared strings, referenced by
# both the header cells A1/B1), and drop the lingering manual cell
# selection (G10) left over from editing, returning the view to the
# default top-left cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "female"
$ws.Range("B1").Value = "male"

# Clear the stale "G10" selection saved in the sheet view by reselecting
# the default top-left cell.
$ws.Range("A1").Select()
